$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 125993.75
$ws.Range("I70").Value = 1090
$ws.Range("K70").Value = 3270
$ws.Range("M70").Value = -3000
$ws.Range("H73").Value = 125993.75
$ws.Range("I73").Value = 1090
$ws.Range("K73").Value = 3270
$ws.Range("M73").Value = -2334
$ws.Range("H96").Value = 2094.25
$ws.Range("I96").Value = 2094.25
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 6282.75
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4909.75
$ws.Range("H98").Value = 3211.2
$ws.Range("H122").Value = 3211.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 878.26666
$ws.Range("I5").Value = 938.9286
$ws.Range("K5").Value = 938.9286
$ws.Range("M5").Value = -826.9286
$ws.Range("H132").Value = 5487.88
$ws.Range("I132").Value = 5487.88
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16463.64
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -13933.64

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 878.26666
$ws.Range("I4").Value = 938.9286
$ws.Range("K4").Value = 938.9286
$ws.Range("M4").Value = -823.9286
$ws.Range("H20").Value = 4168.2
$ws.Range("I20").Value = 3710.75
$ws.Range("K20").Value = 3710.75
$ws.Range("M20").Value = -3463.75
$ws.Range("H58").Value = 35903
$ws.Range("I58").Value = 27709
$ws.Range("J58").Value = 40000
$ws.Range("K58").Value = 27709
$ws.Range("L58").Value = 40000
$ws.Range("M58").Value = -27415
$ws.Range("N58").Value = -40588
$ws.Range("H60").Value = 43697.8
$ws.Range("I60").Value = 20709
$ws.Range("J60").Value = 49445
$ws.Range("K60").Value = 20709
$ws.Range("L60").Value = 49445
$ws.Range("M60").Value = -20110
$ws.Range("N60").Value = -50643

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 29303.5
$ws.Range("I16").Value = 21287.928
$ws.Range("K16").Value = 21287.928
$ws.Range("M16").Value = -21000.928
$ws.Range("H22").Value = 746.7857
$ws.Range("I22").Value = 741.2727
$ws.Range("K22").Value = 741.2727
$ws.Range("M22").Value = -391.2727
$ws.Range("H31").Value = 4494.154
$ws.Range("I31").Value = 3193.4614
$ws.Range("J31").Value = 5794.846
$ws.Range("K31").Value = 3193.4614
$ws.Range("L31").Value = 5794.846
$ws.Range("M31").Value = -2898.4614
$ws.Range("N31").Value = -6384.846
$ws.Range("H34").Value = 4494.154
$ws.Range("I34").Value = 3193.4614
$ws.Range("J34").Value = 5794.846
$ws.Range("K34").Value = 3193.4614
$ws.Range("L34").Value = 5794.846
$ws.Range("M34").Value = -2991.4614
$ws.Range("N34").Value = -6198.846
$ws.Range("H58").Value = 4373.8057
$ws.Range("I58").Value = 4466.5
$ws.Range("J58").Value = 3910.3333
$ws.Range("K58").Value = 4466.5
$ws.Range("L58").Value = 3910.3333
$ws.Range("M58").Value = -4263.5
$ws.Range("N58").Value = -4316.3333
$ws.Range("H99").Value = 4026.1667
$ws.Range("I99").Value = 2724.5715
$ws.Range("K99").Value = 2724.5715
$ws.Range("M99").Value = -1226.5715
$ws.Range("H105").Value = 1370.75
$ws.Range("I105").Value = 1313.5454
$ws.Range("K105").Value = 1313.5454
$ws.Range("M105").Value = 433.4546
$ws.Range("H107").Value = 1430.3462
$ws.Range("I107").Value = 1190.1765
$ws.Range("J107").Value = 1884
$ws.Range("K107").Value = 1190.1765
$ws.Range("L107").Value = 1884
$ws.Range("M107").Value = 729.8235
$ws.Range("N107").Value = -5724
$ws.Range("H113").Value = 29303.5
$ws.Range("I113").Value = 21287.928
$ws.Range("K113").Value = 21287.928
$ws.Range("M113").Value = -19117.928
$ws.Range("H126").Value = 4026.1667
$ws.Range("I126").Value = 2724.5715
$ws.Range("K126").Value = 8173.7145
$ws.Range("M126").Value = -5703.7145
$ws.Range("H132").Value = 3972.1853
$ws.Range("I132").Value = 2970.5908
$ws.Range("K132").Value = 8911.7724
$ws.Range("M132").Value = -6381.7724
$ws.Range("H136").Value = 4373.8057
$ws.Range("I136").Value = 4466.5
$ws.Range("J136").Value = 3910.3333
$ws.Range("K136").Value = 13399.5
$ws.Range("L136").Value = 11730.9999
$ws.Range("M136").Value = -10849.5
$ws.Range("N136").Value = -16830.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1927.5
$ws.Range("I113").Value = 1117
$ws.Range("J113").Value = 2467.8333
$ws.Range("K113").Value = 3351
$ws.Range("L113").Value = 7403.499899999999
$ws.Range("M113").Value = -1181
$ws.Range("N113").Value = -11743.4999
$ws.Range("H131").Value = 2115.6538
$ws.Range("I131").Value = 1164.625
$ws.Range("K131").Value = 3493.875
$ws.Range("M131").Value = 1546.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 245.19048
$ws.Range("J97").Value = 265.66666
$ws.Range("L97").Value = 265.66666
$ws.Range("N97").Value = -1257.66666
$ws.Range("H102").Value = 2126.111
$ws.Range("I102").Value = 2126.111
$ws.Range("K102").Value = 2126.111
$ws.Range("M102").Value = -504.1109999999999
$ws.Range("H124").Value = 94775.5
$ws.Range("J124").Value = 94775.5
$ws.Range("L124").Value = 94775.5
$ws.Range("N124").Value = -104595.5
$ws.Range("H132").Value = 2423
$ws.Range("I132").Value = 2423
$ws.Range("K132").Value = 7269
$ws.Range("M132").Value = -4739

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12818.429
$ws.Range("I40").Value = 8598.1
$ws.Range("K40").Value = 8598.1
$ws.Range("M40").Value = -8462.1
$ws.Range("H55").Value = 192.26666
$ws.Range("I55").Value = 153.44444
$ws.Range("J55").Value = 250.5
$ws.Range("K55").Value = 153.44444
$ws.Range("L55").Value = 250.5
$ws.Range("M55").Value = 19.55556000000001
$ws.Range("N55").Value = -596.5
$ws.Range("H122").Value = 4255.7334
$ws.Range("I122").Value = 4023.4167
$ws.Range("K122").Value = 12070.2501
$ws.Range("M122").Value = -9620.250100000001
$ws.Range("H130").Value = 103821.25
$ws.Range("J130").Value = 103821.25
$ws.Range("L130").Value = 103821.25
$ws.Range("N130").Value = -113861.25
$ws.Range("H132").Value = 19735.17
$ws.Range("I132").Value = 26715.379
$ws.Range("K132").Value = 80146.137
$ws.Range("M132").Value = -77616.137
$ws.Range("H136").Value = 2263.5833
$ws.Range("I136").Value = 2263.5833
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6790.749899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -4240.749899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 883.2727
$ws.Range("I100").Value = 923.6667
$ws.Range("K100").Value = 1847.3334
$ws.Range("M100").Value = -1306.3334
$ws.Range("H125").Value = 70715
$ws.Range("J125").Value = 70715
$ws.Range("L125").Value = 70715
$ws.Range("N125").Value = -80555
$ws.Range("H126").Value = 3473.6562
$ws.Range("I126").Value = 3751.8333
$ws.Range("J126").Value = 2639.125
$ws.Range("K126").Value = 11255.4999
$ws.Range("L126").Value = 7917.375
$ws.Range("M126").Value = -8785.499899999999
$ws.Range("N126").Value = -12857.375
$ws.Range("H132").Value = 2361
$ws.Range("I132").Value = 2361
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7083
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -4553
$ws.Range("H136").Value = 3055.6428
$ws.Range("I136").Value = 3408.875
$ws.Range("K136").Value = 10226.625
$ws.Range("M136").Value = -7676.625
